# Actualización automática 2025-06-16 13:01:14
#
# "VENTAS POR GRUPO" sheet gains a new product-category column ("GRANITO")
# inserted right before the existing "GRIFERIAS" column (old column F),
# shifting GRIFERIAS..SAL SOLUBLE one column to the right, and three more
# new category columns ("NO RESURTIBLES", "PANELES PVC", "PANELES PU")
# appended after the old last column ("SAL SOLUBLE").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "GRANITO" column before the old column F (GRIFERIAS) ---
# This shifts the existing F..N columns (GRIFERIAS..SAL SOLUBLE) to G..O,
# carrying their values/styles/widths along automatically.
$ws.Columns.Item(6).Insert()

$ws.Range("F1").Value2 = "GRANITO"
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 6).Value2 = 0
}
$ws.Range("F19").Value2 = "0 de 17"

# Restore the intended width for the newly inserted column (raw width 13).
$ws.Columns.Item(6).ColumnWidth = 12.17

# --- Append three new columns after the old last column (now column O) ---
# Clone formatting from column O (header bold/center, currency data style,
# centered footer style) onto the new P:R columns before filling values.
$ws.Range("O1:O19").Copy()
$ws.Range("P1:R19").PasteSpecial(-4122)

$ws.Range("P1").Value2 = "NO RESURTIBLES"
$ws.Range("Q1").Value2 = "PANELES PVC"
$ws.Range("R1").Value2 = "PANELES PU"

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 16).Value2 = 0
    $ws.Cells.Item($r, 17).Value2 = 0
    $ws.Cells.Item($r, 18).Value2 = 0
}

$ws.Range("P19").Value2 = "0 de 17"
$ws.Range("Q19").Value2 = "0 de 17"
$ws.Range("R19").Value2 = "0 de 17"

# Target raw column widths: P=20, Q=17, R=16.
$ws.Columns.Item(16).ColumnWidth = 19.17
$ws.Columns.Item(17).ColumnWidth = 16.17
$ws.Columns.Item(18).ColumnWidth = 15.17
